# Fill in the "Need Finding 2" results tables (slide 9 and slide 10).
$p = $ppt.ActivePresentation

# --- Slide 9: two small "Overall Results" tables (one row each, participant "Sean") ---
$s9 = $p.Slides.Item(9)

# Table 5 (shape index 2) -> "Need Finding 1 Colors" summary
$tbl1 = $s9.Shapes.Item(2).Table
$tbl1.Cell(2,2).Shape.TextFrame.TextRange.Text = "28"
$tbl1.Cell(2,3).Shape.TextFrame.TextRange.Text = "0.857 +/- 0.356"
$tbl1.Cell(2,4).Shape.TextFrame.TextRange.Text = "1.57 +/- 0.52"

# Table 6 (shape index 3) -> "Need Finding 2 Colors" summary
$tbl2 = $s9.Shapes.Item(3).Table
$tbl2.Cell(2,2).Shape.TextFrame.TextRange.Text = "168"
$tbl2.Cell(2,3).Shape.TextFrame.TextRange.Text = "0.964 +/- 0.186 "
$tbl2.Cell(2,4).Shape.TextFrame.TextRange.Text = "1.57 +/- 0.43"

# --- Slide 10: "By Color Results" table (one row per color) ---
$s10 = $p.Slides.Item(10)
$tbl = $s10.Shapes.Item(2).Table

# Row 2: Red
$tbl.Cell(2,2).Shape.TextFrame.TextRange.Text = "40"
$tbl.Cell(2,3).Shape.TextFrame.TextRange.Text = "1.000 +/- 0.000"
$tbl.Cell(2,4).Shape.TextFrame.TextRange.Text = "1.46 +/- 0.38"
$tbl.Cell(2,5).Shape.TextFrame.TextRange.Text = "1.000 +/- 0.000"
$tbl.Cell(2,6).Shape.TextFrame.TextRange.Text = "1.927 +/- 0.364"

# Row 3: Brown
$tbl.Cell(3,2).Shape.TextFrame.TextRange.Text = "40"
$tbl.Cell(3,3).Shape.TextFrame.TextRange.Text = "0.500 +/- 0.577"
$tbl.Cell(3,4).Shape.TextFrame.TextRange.Text = "2.12 +/- 0.55"
$tbl.Cell(3,5).Shape.TextFrame.TextRange.Text = "0.917 +/- 0.282"
$tbl.Cell(3,6).Shape.TextFrame.TextRange.Text = "1.627 +/- 0.336"

# Row 4: Green
$tbl.Cell(4,2).Shape.TextFrame.TextRange.Text = "40"
$tbl.Cell(4,3).Shape.TextFrame.TextRange.Text = "0.750 +/- 0.500"
$tbl.Cell(4,4).Shape.TextFrame.TextRange.Text = "1.51 +/- 0.38"
$tbl.Cell(4,5).Shape.TextFrame.TextRange.Text = "0.875 +/- 0.338"
$tbl.Cell(4,6).Shape.TextFrame.TextRange.Text = "1.801 +/- 0.506"

# Row 5: Orange
$tbl.Cell(5,2).Shape.TextFrame.TextRange.Text = "40"
$tbl.Cell(5,3).Shape.TextFrame.TextRange.Text = "0.750 +/- 0.500"
$tbl.Cell(5,4).Shape.TextFrame.TextRange.Text = "1.55 +/- 0.30"
$tbl.Cell(5,5).Shape.TextFrame.TextRange.Text = "0.958 +/- 0.204"
$tbl.Cell(5,6).Shape.TextFrame.TextRange.Text = "1.552 +/- 0.363"

# Row 6: Yellow
$tbl.Cell(6,2).Shape.TextFrame.TextRange.Text = "40"
$tbl.Cell(6,3).Shape.TextFrame.TextRange.Text = "1.000 +/- 0.000"
$tbl.Cell(6,4).Shape.TextFrame.TextRange.Text = "1.12 +/- 0.04"
$tbl.Cell(6,5).Shape.TextFrame.TextRange.Text = "1.000 +/- 0.000"
$tbl.Cell(6,6).Shape.TextFrame.TextRange.Text = "1.223 +/- 0.355"

# Row 7: Blue
$tbl.Cell(7,2).Shape.TextFrame.TextRange.Text = "40"
$tbl.Cell(7,3).Shape.TextFrame.TextRange.Text = "1.000 +/- 0.000"
$tbl.Cell(7,4).Shape.TextFrame.TextRange.Text = "1.67 +/- 0.99"
$tbl.Cell(7,5).Shape.TextFrame.TextRange.Text = "1.000 +/- 0.000"
$tbl.Cell(7,6).Shape.TextFrame.TextRange.Text = "1.366 +/- 0.364"

# Row 8: Purple
$tbl.Cell(8,2).Shape.TextFrame.TextRange.Text = "40"
$tbl.Cell(8,3).Shape.TextFrame.TextRange.Text = "1.000 +/- 0.000"
$tbl.Cell(8,4).Shape.TextFrame.TextRange.Text = "1.55 +/- 0.29"
$tbl.Cell(8,5).Shape.TextFrame.TextRange.Text = "1.000 +/- 0.000"
$tbl.Cell(8,6).Shape.TextFrame.TextRange.Text = "1.500 +/- 0.331"
